$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5882880
$ws.Range("I33").Value = 8333992.5
$ws.Range("J33").Value = 210
$ws.Range("K33").Value = 8333992.5
$ws.Range("L33").Value = 210
$ws.Range("M33").Value = -8333763.5
$ws.Range("N33").Value = -668
$ws.Range("H64").Value = 6357.2856
$ws.Range("I64").Value = 6145.3335
$ws.Range("J64").Value = 6738.8
$ws.Range("K64").Value = 6145.3335
$ws.Range("L64").Value = 6738.8
$ws.Range("M64").Value = -5897.3335
$ws.Range("N64").Value = -7234.8
$ws.Range("H67").Value = 6357.2856
$ws.Range("I67").Value = 6145.3335
$ws.Range("J67").Value = 6738.8
$ws.Range("K67").Value = 6145.3335
$ws.Range("L67").Value = 6738.8
$ws.Range("M67").Value = -5287.3335
$ws.Range("N67").Value = -8454.799999999999
$ws.Range("H70").Value = 7084.4443
$ws.Range("I70").Value = 5663.6665
$ws.Range("J70").Value = 7794.8335
$ws.Range("K70").Value = 16990.9995
$ws.Range("L70").Value = 23384.5005
$ws.Range("M70").Value = -16720.9995
$ws.Range("N70").Value = -23924.5005
$ws.Range("H73").Value = 7084.4443
$ws.Range("I73").Value = 5663.6665
$ws.Range("J73").Value = 7794.8335
$ws.Range("K73").Value = 16990.9995
$ws.Range("L73").Value = 23384.5005
$ws.Range("M73").Value = -16054.9995
$ws.Range("N73").Value = -25256.5005
$ws.Range("H137").Value = 13897217
$ws.Range("I137").Value = 31250932
$ws.Range("J137").Value = 14244.8
$ws.Range("K137").Value = 93752796
$ws.Range("L137").Value = 42734.39999999999
$ws.Range("M137").Value = -93750246
$ws.Range("N137").Value = -47834.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 43729.043
$ws.Range("I32").Value = 43729.043
$ws.Range("K32").Value = 43729.043
$ws.Range("M32").Value = -43442.043
$ws.Range("H45").Value = 2446.0833
$ws.Range("J45").Value = 2983
$ws.Range("L45").Value = 2983
$ws.Range("N45").Value = -3737
$ws.Range("H132").Value = 3336.2068
$ws.Range("I132").Value = 2335.7083
$ws.Range("K132").Value = 7007.124899999999
$ws.Range("M132").Value = -4477.124899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 206.5
$ws.Range("I11").Value = 206.5
$ws.Range("K11").Value = 206.5
$ws.Range("M11").Value = -66.5
$ws.Range("H105").Value = 52645836
$ws.Range("I105").Value = 52645836
$ws.Range("K105").Value = 52645836
$ws.Range("M105").Value = -52644089
$ws.Range("H110").Value = 45702
$ws.Range("J110").Value = 45702
$ws.Range("L110").Value = 45702
$ws.Range("N110").Value = -53882
$ws.Range("H130").Value = 119966.664
$ws.Range("J130").Value = 119966.664
$ws.Range("L130").Value = 119966.664
$ws.Range("N130").Value = -130006.664
$ws.Range("H134").Value = 2995.5925
$ws.Range("I134").Value = 1437.2858
$ws.Range("J134").Value = 8449.666999999999
$ws.Range("K134").Value = 4311.857400000001
$ws.Range("L134").Value = 25349.001
$ws.Range("M134").Value = -1776.857400000001
$ws.Range("N134").Value = -30419.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 29996.334
$ws.Range("J62").Value = 29996.334
$ws.Range("L62").Value = 29996.334
$ws.Range("N62").Value = -31244.334
$ws.Range("H65").Value = 29996.334
$ws.Range("J65").Value = 29996.334
$ws.Range("L65").Value = 149981.67
$ws.Range("N65").Value = -156221.67
$ws.Range("H74").Value = 80000
$ws.Range("J74").Value = 80000
$ws.Range("L74").Value = 80000
$ws.Range("N74").Value = -81748
$ws.Range("H77").Value = 80000
$ws.Range("J77").Value = 80000
$ws.Range("L77").Value = 240000
$ws.Range("N77").Value = -248736
$ws.Range("H132").Value = 158215.75
$ws.Range("I132").Value = 7240.5
$ws.Range("J132").Value = 309191
$ws.Range("K132").Value = 21721.5
$ws.Range("L132").Value = 927573
$ws.Range("M132").Value = -19191.5
$ws.Range("N132").Value = -932633
$ws.Range("H134").Value = 13498.615
$ws.Range("I134").Value = 13448.454
$ws.Range("J134").Value = 13774.5
$ws.Range("K134").Value = 40345.362
$ws.Range("L134").Value = 41323.5
$ws.Range("M134").Value = -37810.362
$ws.Range("N134").Value = -46393.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4914.727
$ws.Range("J64").Value = 3009.125
$ws.Range("L64").Value = 9027.375
$ws.Range("N64").Value = -9567.375
$ws.Range("H67").Value = 4914.727
$ws.Range("J67").Value = 3009.125
$ws.Range("L67").Value = 9027.375
$ws.Range("N67").Value = -10899.375
$ws.Range("H74").Value = 9128.25
$ws.Range("I74").Value = 9004.333000000001
$ws.Range("K74").Value = 27012.999
$ws.Range("M74").Value = -25951.999
$ws.Range("H75").Value = 1512.5
$ws.Range("J75").Value = 1709.8572
$ws.Range("L75").Value = 5129.571599999999
$ws.Range("N75").Value = -7125.571599999999
$ws.Range("H77").Value = 9128.25
$ws.Range("I77").Value = 9004.333000000001
$ws.Range("K77").Value = 81038.997
$ws.Range("M77").Value = -75734.997
$ws.Range("H78").Value = 1512.5
$ws.Range("J78").Value = 1709.8572
$ws.Range("L78").Value = 15388.7148
$ws.Range("N78").Value = -25372.7148
$ws.Range("H122").Value = 11689.857
$ws.Range("J122").Value = 2000.5
$ws.Range("L122").Value = 18004.5
$ws.Range("N122").Value = -22904.5
$ws.Range("H128").Value = 436622
$ws.Range("I128").Value = 436622
$ws.Range("K128").Value = 1309866
$ws.Range("M128").Value = -1304886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2799.3333
$ws.Range("I80").Value = 1632.3334
$ws.Range("K80").Value = 1632.3334
$ws.Range("M80").Value = -634.3334
$ws.Range("H83").Value = 2799.3333
$ws.Range("I83").Value = 1632.3334
$ws.Range("K83").Value = 8161.666999999999
$ws.Range("M83").Value = -3169.666999999999
$ws.Range("H113").Value = 28097.154
$ws.Range("I113").Value = 1546.6666
$ws.Range("K113").Value = 1546.6666
$ws.Range("M113").Value = 623.3334
$ws.Range("H132").Value = 4277.242
$ws.Range("I132").Value = 2109.611
$ws.Range("K132").Value = 6328.833
$ws.Range("M132").Value = -3798.833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2750.6667
$ws.Range("I22").Value = 1689.0769
$ws.Range("J22").Value = 4203.3687
$ws.Range("K22").Value = 1689.0769
$ws.Range("L22").Value = 4203.3687
$ws.Range("M22").Value = -1394.0769
$ws.Range("N22").Value = -4793.3687
$ws.Range("H27").Value = 2750.6667
$ws.Range("I27").Value = 1689.0769
$ws.Range("J27").Value = 4203.3687
$ws.Range("K27").Value = 1689.0769
$ws.Range("L27").Value = 4203.3687
$ws.Range("M27").Value = -1582.0769
$ws.Range("N27").Value = -4417.3687
$ws.Range("H122").Value = 3952.5715
$ws.Range("I122").Value = 3952.5715
$ws.Range("K122").Value = 11857.7145
$ws.Range("M122").Value = -9407.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5454.091
$ws.Range("I132").Value = 3018.625
$ws.Range("J132").Value = 7746.294
$ws.Range("K132").Value = 9055.875
$ws.Range("L132").Value = 23238.882
$ws.Range("M132").Value = -6525.875
$ws.Range("N132").Value = -28298.882
